$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) CGPA cell: "8.6" -> "8.66"
# ---------------------------------------------------------------------------
$cgpaRange = $d.Content
$found = $cgpaRange.Find.Execute("8.6", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
if ($found) {
    $insertPos = $cgpaRange.End
    $ins = $d.Range($insertPos, $insertPos)
    $ins.InsertAfter("6")
}

# ---------------------------------------------------------------------------
# 2) Academic participation bullet list: insert a new "NPTEL ..." bullet
#    before the existing "Cloud quiz ..." bullet, and push the existing
#    "General coding quiz" bullet into a brand-new paragraph so the final
#    order is: NPTEL ... / Cloud quiz ... / General coding quiz
# ---------------------------------------------------------------------------

# Locate the "Cloud quiz organized..." paragraph (currently 1st of the two
# bullets under ACADEMIC PARTICIPATION).
$cloudIndex = -1
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Cloud quiz organized by CloudZone club*") {
        $cloudIndex = $i
    }
}

if ($cloudIndex -gt 0) {
    # --- Turn this paragraph's text into the new NPTEL bullet (in place,
    #     keeping the paragraph / numbering formatting untouched). ---
    $pCloud = $d.Paragraphs.Item($cloudIndex)
    $rCloud = $pCloud.Range
    $rCloudNoMark = $d.Range($rCloud.Start, $rCloud.End - 1)
    $rCloudNoMark.Text = "NPTEL Introduction to Artificial Intelligence Certificate which had only 2.49% success rate"

    # --- The next paragraph is "General coding quiz"; turn its text into
    #     the "Cloud quiz ... position" bullet (with superscript "st"). ---
    $genIndex = $cloudIndex + 1
    $pGen = $d.Paragraphs.Item($genIndex)
    $rGen = $pGen.Range
    $rGenNoMark = $d.Range($rGen.Start, $rGen.End - 1)
    $rGenNoMark.Text = "Cloud quiz organized by CloudZone club " + [char]0x2013 + " 1st position"

    $pGen2 = $d.Paragraphs.Item($genIndex)
    $genText = $pGen2.Range.Text
    $stOffset = $genText.IndexOf("st")
    $stStart = $rGen.Start + $stOffset
    $stRange = $d.Range($stStart, $stStart + 2)
    $stRange.Font.Superscript = $true

    # --- Append a brand-new paragraph after it holding "General coding
    #     quiz" (this is the paragraph that picks up the bullet formatting
    #     automatically from its predecessor). ---
    $pGen3 = $d.Paragraphs.Item($genIndex)
    $pGen3.Range.InsertParagraphAfter()
    $pNew = $d.Paragraphs.Item($genIndex + 1)
    $pNew.Range.InsertAfter("General coding quiz")
}
